# Auto-generated edit script: updates numeric leve-profit values
# per the commit diff (scheduled runner refresh of market-price-derived columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 63
  $ws.Range("H63").Value = 26544.908
  $ws.Range("J63").Value = 26544.908
  $ws.Range("L63").Value = 26544.908
  $ws.Range("N63").Value = -27792.908
  # Row 66
  $ws.Range("H66").Value = 26544.908
  $ws.Range("J66").Value = 26544.908
  $ws.Range("L66").Value = 79634.724
  $ws.Range("N66").Value = -85874.724
  # Row 75
  $ws.Range("H75").Value = 40666.668
  $ws.Range("J75").Value = 40666.668
  $ws.Range("L75").Value = 40666.668
  $ws.Range("N75").Value = -42538.668
  # Row 78
  $ws.Range("H78").Value = 40666.668
  $ws.Range("J78").Value = 40666.668
  $ws.Range("L78").Value = 122000.004
  $ws.Range("N78").Value = -131360.004
  # Row 93
  $ws.Range("H93").Value = 42500
  $ws.Range("J93").Value = 42500
  $ws.Range("L93").Value = 42500
  $ws.Range("N93").Value = -47492
  # Row 108
  $ws.Range("H108").Value = 38538.75
  $ws.Range("J108").Value = 38538.75
  $ws.Range("L108").Value = 38538.75
  $ws.Range("N108").Value = -46218.75
  # Row 118
  $ws.Range("H118").Value = 149376.03
  $ws.Range("I118").Value = 667160.8
  $ws.Range("J118").Value = 1437.5238
  $ws.Range("K118").Value = 2001482.4
  $ws.Range("L118").Value = 4312.5714
  $ws.Range("M118").Value = -1999825.4
  $ws.Range("N118").Value = -7626.5714
  # Row 130
  $ws.Range("H130").Value = 36663.332
  $ws.Range("J130").Value = 36663.332
  $ws.Range("L130").Value = 36663.332
  $ws.Range("N130").Value = -46703.332
  # Row 132
  $ws.Range("H132").Value = 3004.2046
  $ws.Range("I132").Value = 2572.6667
  $ws.Range("J132").Value = 3928.9285
  $ws.Range("K132").Value = 7718.000100000001
  $ws.Range("L132").Value = 11786.7855
  $ws.Range("M132").Value = -5188.000100000001
  $ws.Range("N132").Value = -16846.7855

$ws = $wb.Worksheets.Item("ARM")
  # Row 32
  $ws.Range("H32").Value = 25454.842
  $ws.Range("I32").Value = 6452.349
  $ws.Range("J32").Value = 117543.84
  $ws.Range("K32").Value = 6452.349
  $ws.Range("L32").Value = 117543.84
  $ws.Range("M32").Value = -6165.349
  $ws.Range("N32").Value = -118117.84
  # Row 44
  $ws.Range("H44").Value = 35044.5
  $ws.Range("J44").Value = 35044.5
  $ws.Range("L44").Value = 35044.5
  $ws.Range("N44").Value = -36020.5
  # Row 103
  $ws.Range("H103").Value = 39678
  $ws.Range("J103").Value = 39678
  $ws.Range("L103").Value = 39678
  $ws.Range("N103").Value = -42022
  # Row 127
  $ws.Range("H127").Value = 0
  $ws.Range("J127").Value = 0
  $ws.Range("L127").Value = 0
  $ws.Range("N127").ClearContents()
  # Row 128
  $ws.Range("H128").Value = 0
  $ws.Range("J128").Value = 0
  $ws.Range("L128").Value = 0
  $ws.Range("N128").ClearContents()
  # Row 134
  $ws.Range("H134").Value = 37493.332
  $ws.Range("J134").Value = 37493.332
  $ws.Range("L134").Value = 37493.332
  $ws.Range("N134").Value = -47633.332
  # Row 135
  $ws.Range("H135").Value = 38510.4
  $ws.Range("J135").Value = 38510.4
  $ws.Range("L135").Value = 38510.4
  $ws.Range("N135").Value = -48650.4

$ws = $wb.Worksheets.Item("BSM")
  # Row 82
  $ws.Range("H82").Value = 59394.562
  $ws.Range("J82").Value = 30769.3
  $ws.Range("L82").Value = 30769.3
  $ws.Range("N82").Value = -31535.3
  # Row 85
  $ws.Range("H85").Value = 59394.562
  $ws.Range("J85").Value = 30769.3
  $ws.Range("L85").Value = 30769.3
  $ws.Range("N85").Value = -33421.3
  # Row 122
  $ws.Range("H122").Value = 41930
  $ws.Range("J122").Value = 41930
  $ws.Range("L122").Value = 41930
  $ws.Range("N122").Value = -51730
  # Row 124
  $ws.Range("H124").Value = 42308
  $ws.Range("J124").Value = 42308
  $ws.Range("L124").Value = 42308
  $ws.Range("N124").Value = -52128
  # Row 132
  $ws.Range("H132").Value = 30918.572
  $ws.Range("J132").Value = 30918.572
  $ws.Range("L132").Value = 30918.572
  $ws.Range("N132").Value = -41038.572
  # Row 135
  $ws.Range("H135").Value = 66773.84
  $ws.Range("J135").Value = 66773.84
  $ws.Range("L135").Value = 66773.84
  $ws.Range("N135").Value = -76913.84

$ws = $wb.Worksheets.Item("CRP")
  # Row 20
  $ws.Range("H20").Value = 48874.75
  $ws.Range("J20").Value = 48874.75
  $ws.Range("L20").Value = 48874.75
  $ws.Range("N20").Value = -49346.75
  # Row 30
  $ws.Range("H30").Value = 48874.75
  $ws.Range("J30").Value = 48874.75
  $ws.Range("L30").Value = 48874.75
  $ws.Range("N30").Value = -49056.75
  # Row 41
  $ws.Range("H41").Value = 15869.286
  $ws.Range("J41").Value = 20637
  $ws.Range("L41").Value = 20637
  $ws.Range("N41").Value = -21493
  # Row 127
  $ws.Range("H127").Value = 54750
  $ws.Range("J127").Value = 54500
  $ws.Range("L127").Value = 54500
  $ws.Range("N127").Value = -64420
  # Row 128
  $ws.Range("H128").Value = 48874.75
  $ws.Range("J128").Value = 48874.75
  $ws.Range("L128").Value = 48874.75
  $ws.Range("N128").Value = -58834.75
  # Row 130
  $ws.Range("H130").Value = 29953.334
  $ws.Range("J130").Value = 29953.334
  $ws.Range("L130").Value = 29953.334
  $ws.Range("N130").Value = -39993.334

$ws = $wb.Worksheets.Item("CUL")
  # Row 92
  $ws.Range("H92").Value = 787.2222
  $ws.Range("I92").Value = 569
  $ws.Range("J92").Value = 1060
  $ws.Range("K92").Value = 1707
  $ws.Range("L92").Value = 3180
  $ws.Range("M92").Value = -459
  $ws.Range("N92").Value = -5676

$ws = $wb.Worksheets.Item("GSM")
  # Row 43
  $ws.Range("H43").Value = 3568.7273
  $ws.Range("J43").Value = 8226.5
  $ws.Range("L43").Value = 8226.5
  $ws.Range("N43").Value = -8528.5
  # Row 93
  $ws.Range("H93").Value = 13793.5
  $ws.Range("J93").Value = 13793.5
  $ws.Range("L93").Value = 13793.5
  $ws.Range("N93").Value = -17537.5
  # Row 127
  $ws.Range("H127").Value = 55993.332
  $ws.Range("J127").Value = 55993.332
  $ws.Range("L127").Value = 55993.332
  $ws.Range("N127").Value = -65913.33199999999
  # Row 128
  $ws.Range("H128").Value = 54100
  $ws.Range("J128").Value = 54100
  $ws.Range("L128").Value = 54100
  $ws.Range("N128").Value = -64060
  # Row 132
  $ws.Range("H132").Value = 2894.6365
  $ws.Range("I132").Value = 2593.6667
  $ws.Range("J132").Value = 4249
  $ws.Range("K132").Value = 7781.000100000001
  $ws.Range("L132").Value = 12747
  $ws.Range("M132").Value = -5251.000100000001
  $ws.Range("N132").Value = -17807
  # Row 133
  $ws.Range("H133").Value = 26424.285
  $ws.Range("J133").Value = 26424.285
  $ws.Range("L133").Value = 26424.285
  $ws.Range("N133").Value = -36544.285
  # Row 135
  $ws.Range("H135").Value = 53712.5
  $ws.Range("J135").Value = 53712.5
  $ws.Range("L135").Value = 53712.5
  $ws.Range("N135").Value = -63852.5

$ws = $wb.Worksheets.Item("LTW")
  # Row 125
  $ws.Range("H125").Value = 35357.5
  $ws.Range("J125").Value = 35357.5
  $ws.Range("L125").Value = 35357.5
  $ws.Range("N125").Value = -45197.5
  # Row 130
  $ws.Range("H130").Value = 54424.5
  $ws.Range("J130").Value = 54424.5
  $ws.Range("L130").Value = 54424.5
  $ws.Range("N130").Value = -64464.5
  # Row 133
  $ws.Range("H133").Value = 56292.5
  $ws.Range("J133").Value = 56292.5
  $ws.Range("L133").Value = 56292.5
  $ws.Range("N133").Value = -61352.5
  # Row 134
  $ws.Range("H134").Value = 70214.5
  $ws.Range("I134").Value = 0
  $ws.Range("J134").Value = 70214.5
  $ws.Range("K134").Value = 0
  $ws.Range("L134").Value = 70214.5
  $ws.Range("M134").ClearContents()
  $ws.Range("N134").Value = -80354.5

$ws = $wb.Worksheets.Item("WVR")
  # Row 129
  $ws.Range("H129").Value = 39322.25
  $ws.Range("J129").Value = 39322.25
  $ws.Range("L129").Value = 39322.25
  $ws.Range("N129").Value = -49322.25
  # Row 135
  $ws.Range("H135").Value = 52844.54
  $ws.Range("J135").Value = 53498.25
  $ws.Range("L135").Value = 53498.25
  $ws.Range("N135").Value = -63638.25
